$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new blank row at position 5. This pushes the old "Source" row
#    (old row 5) down to row 6, while the old row 4 (label + numbers) stays
#    at row 4 (it will be fully re-populated with new label/values below).
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Insert()

# ---------------------------------------------------------------------------
# 2. Normalise the numeric-cell formatting on row 4 (B4:I4) so that none of
#    them carry a horizontal alignment override (the new layout has no
#    alignment set on the data cells). Base this off B4's existing format.
# ---------------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("B4:I4").PasteSpecial(-4122)

# Remove the bottom border from A4 (keep the top border).
$ws.Range("A4").Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------------
# 3. Populate the newly created row 5 using row 4 as a formatting template.
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
# A5 needs bottom border instead of top border.
$ws.Range("A5").Borders.Item(8).LineStyle = -4142
$ws.Range("A5").Borders.Item(9).LineStyle = 1

$ws.Range("B4:I4").Copy()
$ws.Range("B5:I5").PasteSpecial(-4122)
# I5 additionally needs a bottom border.
$ws.Range("I5").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# 4. Row 6 (previously row 5, the merged "Source" row) - remove the top
#    border that used to sit on the label cell A6.
# ---------------------------------------------------------------------------
$ws.Range("A6").Borders.Item(8).LineStyle = -4142

# ---------------------------------------------------------------------------
# 5. New row 4 content: "family with disabilities Persons "
# ---------------------------------------------------------------------------
$ws.Range("A4").Value2 = "family with disabilities Persons "
$ws.Range("B4").Value2 = 549
$ws.Range("C4").Value2 = 543
$ws.Range("D4").Value2 = 538
$ws.Range("E4").Value2 = 598
$ws.Range("F4").Value2 = 644
$ws.Range("G4").Value2 = 692
$ws.Range("H4").Value2 = 756
$ws.Range("I4").Value2 = 926

# ---------------------------------------------------------------------------
# 6. New row 5 content: "disabilities Persons "
# ---------------------------------------------------------------------------
$ws.Range("A5").Value2 = "disabilities Persons "
$ws.Range("B5").Value2 = 605
$ws.Range("C5").Value2 = 603
$ws.Range("D5").Value2 = 595
$ws.Range("E5").Value2 = 654
$ws.Range("F5").Value2 = 698
$ws.Range("G5").Value2 = 756
$ws.Range("H5").Value2 = 827
$ws.Range("I5").Value2 = 1012

# ---------------------------------------------------------------------------
# 7. Row 1 title - new wording, centred + wrapped, merged across A1:I1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Marneuli Municipality"
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").WrapText = $true
$ws.Range("A1:I1").Merge()
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------------
# 8. Row 3, cell A3: font changes from Arial to Sylfaen (size 11).
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Family = 1
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# 9. Row heights.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.5
$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 21
$ws.Rows.Item(6).RowHeight = 27.75

Write-Output "structure done"
